$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quality_comparison")
$cell = $ws.Cells.Item(1,3)
$cell.Font.Bold = $true
Write-Host "Bold after: " $cell.Font.Bold
